$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 2 (a blank row except for an "1 кв" note in column F) is
# removed entirely; everything below it shifts up by one row.
$ws.Rows("2").Delete()

# The removed sub-header text is folded into the year header in F1 as a
# second line ("2023" + newline + "1 кв").
$ws.Range("F1").Value = "2023`n1 кв"

# Let Excel recompute the row height for the now two-line header instead
# of leaving an explicit custom height behind.
$ws.Rows("1").AutoFit()

# The average-wage figure for 2023 (now row 28 after the shift) is revised.
$ws.Range("F28").Value = 345.3
